$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2750
$ws.Range("J48").Value = 2750
$ws.Range("L48").Value = 8250
$ws.Range("N48").Value = -8834
$ws.Range("H56").Value = 2750
$ws.Range("J56").Value = 2750
$ws.Range("L56").Value = 8250
$ws.Range("N56").Value = -9318
$ws.Range("H58").Value = 2917.25
$ws.Range("I58").Value = 241.33333
$ws.Range("J58").Value = 4522.8
$ws.Range("K58").Value = 723.99999
$ws.Range("L58").Value = 13568.4
$ws.Range("M58").Value = -573.99999
$ws.Range("N58").Value = -13868.4
$ws.Range("H138").Value = 4359.815
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4359.815
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 13079.445
$ws.Range("N138").Value = -23359.445
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 900000
$ws.Range("J62").Value = 900000
$ws.Range("L62").Value = 900000
$ws.Range("N62").Value = -901248
$ws.Range("H65").Value = 900000
$ws.Range("J65").Value = 900000
$ws.Range("L65").Value = 2700000
$ws.Range("N65").Value = -2706240
$ws.Range("H74").Value = 3141.923
$ws.Range("I74").Value = 2731.4546
$ws.Range("J74").Value = 5399.5
$ws.Range("K74").Value = 2731.4546
$ws.Range("L74").Value = 5399.5
$ws.Range("M74").Value = -1857.4546
$ws.Range("N74").Value = -7147.5
$ws.Range("H77").Value = 3141.923
$ws.Range("I77").Value = 2731.4546
$ws.Range("J77").Value = 5399.5
$ws.Range("K77").Value = 13657.273
$ws.Range("L77").Value = 26997.5
$ws.Range("M77").Value = -9289.273000000001
$ws.Range("N77").Value = -35733.5
$ws.Range("H88").Value = 9272915
$ws.Range("J88").Value = 27779582
$ws.Range("L88").Value = 27779582
$ws.Range("N88").Value = -27780394
$ws.Range("H91").Value = 9272915
$ws.Range("J91").Value = 27779582
$ws.Range("L91").Value = 27779582
$ws.Range("N91").Value = -27782390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3662.2878
$ws.Range("I31").Value = 2955.875
$ws.Range("J31").Value = 4065.9524
$ws.Range("K31").Value = 2955.875
$ws.Range("L31").Value = 4065.9524
$ws.Range("M31").Value = -2660.875
$ws.Range("N31").Value = -4655.9524
$ws.Range("H34").Value = 3662.2878
$ws.Range("I34").Value = 2955.875
$ws.Range("J34").Value = 4065.9524
$ws.Range("K34").Value = 2955.875
$ws.Range("L34").Value = 4065.9524
$ws.Range("M34").Value = -2753.875
$ws.Range("N34").Value = -4469.9524
$ws.Range("H105").Value = 493.375
$ws.Range("I105").Value = 517.25
$ws.Range("J105").Value = 469.5
$ws.Range("K105").Value = 517.25
$ws.Range("L105").Value = 469.5
$ws.Range("M105").Value = 1229.75
$ws.Range("N105").Value = -3963.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 272.5
$ws.Range("I12").Value = 266
$ws.Range("J12").Value = 273.42856
$ws.Range("K12").Value = 798
$ws.Range("L12").Value = 820.28568
$ws.Range("M12").Value = -625
$ws.Range("N12").Value = -1166.28568
$ws.Range("H18").Value = 1156.7142
$ws.Range("I18").Value = 226.33333
$ws.Range("K18").Value = 678.99999
$ws.Range("M18").Value = -509.99999
$ws.Range("H23").Value = 3560.3635
$ws.Range("J23").Value = 2724.4285
$ws.Range("L23").Value = 8173.2855
$ws.Range("N23").Value = -8643.2855
$ws.Range("H36").Value = 5328.1665
$ws.Range("J36").Value = 5994
$ws.Range("L36").Value = 17982
$ws.Range("N36").Value = -18320
$ws.Range("H114").Value = 1619.5
$ws.Range("J114").Value = 1773.4
$ws.Range("L114").Value = 5320.200000000001
$ws.Range("N114").Value = -11828.2
$ws.Range("H117").Value = 2285.4211
$ws.Range("I117").Value = 1128.8334
$ws.Range("J117").Value = 2819.2307
$ws.Range("K117").Value = 3386.5002
$ws.Range("L117").Value = 8457.6921
$ws.Range("M117").Value = 55.49980000000005
$ws.Range("N117").Value = -15341.6921
$ws.Range("H124").Value = 1000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 1000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 3000
$ws.Range("N124").Value = -12820
$ws.Range("H129").Value = 1680.75
$ws.Range("I129").Value = 1165
$ws.Range("J129").Value = 1901.7858
$ws.Range("K129").Value = 3495
$ws.Range("L129").Value = 5705.357400000001
$ws.Range("M129").Value = 1505
$ws.Range("N129").Value = -15705.3574
$ws.Range("M124").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 60476.145
$ws.Range("J135").Value = 60476.145
$ws.Range("L135").Value = 60476.145
$ws.Range("N135").Value = -70616.14499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 68001
$ws.Range("I100").Value = 86832.08
$ws.Range("K100").Value = 86832.08
$ws.Range("M100").Value = -86291.08
$ws.Range("H108").Value = 1000000
$ws.Range("J108").Value = 1000000
$ws.Range("L108").Value = 1000000
$ws.Range("N108").Value = -1007680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 32495.9
$ws.Range("I49").Value = 24998
$ws.Range("J49").Value = 33329
$ws.Range("K49").Value = 24998
$ws.Range("L49").Value = 33329
$ws.Range("M49").Value = -24768
$ws.Range("N49").Value = -33789
$ws.Range("H62").Value = 13892008
$ws.Range("I62").Value = 2750
$ws.Range("J62").Value = 22225562
$ws.Range("K62").Value = 2750
$ws.Range("L62").Value = 22225562
$ws.Range("M62").Value = -2126
$ws.Range("N62").Value = -22226810
$ws.Range("H65").Value = 13892008
$ws.Range("I65").Value = 2750
$ws.Range("J65").Value = 22225562
$ws.Range("K65").Value = 13750
$ws.Range("L65").Value = 111127810
$ws.Range("M65").Value = -10630
$ws.Range("N65").Value = -111134050
$ws.Range("H81").Value = 66672130
$ws.Range("I81").Value = 5000
$ws.Range("K81").Value = 10000
$ws.Range("M81").Value = -8939
$ws.Range("H84").Value = 66672130
$ws.Range("I84").Value = 5000
$ws.Range("K84").Value = 50000
$ws.Range("M84").Value = -44696
$ws.Range("H132").Value = 3322.9565
$ws.Range("I132").Value = 1827.5385
$ws.Range("K132").Value = 5482.6155
$ws.Range("M132").Value = -2952.6155
